$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the style from the existing header cell H1 so I1/J1 match the
# bold/centered/bordered header style used by the rest of row 1.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Data rows 2-32 for columns I (I0) and J (IF)
$data = @{
    2  = @(7, 7)
    3  = @(1, 6)
    4  = @(1, 1)
    5  = @(1, 6)
    6  = @(1, 6)
    7  = @(1, 6)
    8  = @(1, 5)
    9  = @(1, 3)
    10 = @(1, 6)
    11 = @(1, 6)
    12 = @(1, 6)
    13 = @(1, 5)
    14 = @(1, 5)
    15 = @(1, 5)
    16 = @(1, 6)
    17 = @(1, 6)
    18 = @(1, 6)
    19 = @(1, 5)
    20 = @(1, 7)
    21 = @(1, 5)
    22 = @(1, 6)
    23 = @(1, 5)
    24 = @(1, 7)
    25 = @(1, 6)
    26 = @(1, 6)
    27 = @(1, 6)
    28 = @(1, 5)
    29 = @(1, 4)
    30 = @(5, 7)
    31 = @(1, 3)
    32 = @(3, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
